$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Append three new rows to the (first) table with the new CMU Id Card tag
#    values.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)

$newRow = $t.Rows.Add()
$i = $newRow.Index
$t.Cell($i, 1).Range.Text = "Smruthi Id Card"
$t.Cell($i, 2).Range.Text = "Andrew id"
$t.Cell($i, 3).Range.Text = "775284314735"

$newRow = $t.Rows.Add()
$i = $newRow.Index
$t.Cell($i, 1).Range.Text = "Sakthi Id Card"
$t.Cell($i, 2).Range.Text = "Andrew id"
$t.Cell($i, 3).Range.Text = "428573490223"

$newRow = $t.Rows.Add()
$i = $newRow.Index
$t.Cell($i, 1).Range.Text = "Rishabh Id Card"
$t.Cell($i, 2).Range.Text = "Andrew id"
$t.Cell($i, 3).Range.Text = "881298603775"

# ---------------------------------------------------------------------------
# 2. Insert a new empty paragraph right after the table, i.e. just before the
#    pre-existing empty paragraph that sits between the table and the
#    "_GoBack" bookmark paragraph.
# ---------------------------------------------------------------------------
$afterTable = $t.Range.End
$insertionPoint = $d.Range($afterTable, $afterTable)
$insertionPoint.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 3. Remove the trailing empty paragraph that used to sit right before the
#    closing <w:sectPr> (after the "_GoBack" bookmark paragraph), while
#    keeping the bookmark itself intact.
#
#    Word will not let the very last paragraph mark of the document be
#    deleted directly (it anchors the section properties), and merging the
#    bookmark paragraph's own mark forward while that paragraph holds no run
#    content drops the bookmark along with it. So: give the bookmark
#    paragraph a tiny run of placeholder text, delete its paragraph mark
#    (which now correctly merges the bookmark + text into the following,
#    final paragraph), and then strip the placeholder text back out again.
# ---------------------------------------------------------------------------
$placeholder = "ZZPLACEHOLDERZZ"

$bookmarkRange = $d.Bookmarks.Item("_GoBack").Range
$bookmarkStart = $bookmarkRange.Start
$marker = $d.Range($bookmarkStart, $bookmarkStart)
$marker.InsertBefore($placeholder)

$bookmarkRange = $d.Bookmarks.Item("_GoBack").Range
$paraMarkPos = $bookmarkRange.Start
$paraMark = $d.Range($paraMarkPos, $paraMarkPos + 1)
$paraMark.Delete()

$placeholderRange = $d.Range($bookmarkStart, $bookmarkStart + $placeholder.Length)
$placeholderRange.Text = ""
